$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B text updates
$ws.Range("B6").Value  = "<day>"
$ws.Range("B9").Value  = "<be>"
$ws.Range("B14").Value = "<be>"
$ws.Range("B18").Value = "<a>"
$ws.Range("B24").Value = "<there>"
$ws.Range("B33").Value = "<by>"
$ws.Range("B36").Value = "<it>"
$ws.Range("B38").Value = "<many>"
$ws.Range("B48").Value = "<alt>"
$ws.Range("B50").Value = "<xackspace>"
$ws.Range("B51").Value = "<and>"

# Column C numeric updates
$ws.Range("C2").Value  = 6
$ws.Range("C3").Value  = 8
$ws.Range("C9").Value  = 10
$ws.Range("C10").Value = 4
$ws.Range("C12").Value = 6
$ws.Range("C13").Value = 10
$ws.Range("C14").Value = 5
$ws.Range("C16").Value = 8
$ws.Range("C17").Value = 10
$ws.Range("C18").Value = 11
$ws.Range("C19").Value = 12
$ws.Range("C20").Value = 6
$ws.Range("C21").Value = 8
$ws.Range("C23").Value = 3
$ws.Range("C24").Value = 8
$ws.Range("C25").Value = 11
$ws.Range("C26").Value = 9
$ws.Range("C28").Value = 10
$ws.Range("C29").Value = 7
$ws.Range("C30").Value = 15
$ws.Range("C31").Value = 6
$ws.Range("C32").Value = 7
$ws.Range("C33").Value = 7
$ws.Range("C34").Value = 6
$ws.Range("C35").Value = 7
$ws.Range("C36").Value = 10
$ws.Range("C37").Value = 12
$ws.Range("C38").Value = 7
$ws.Range("C39").Value = 5
$ws.Range("C40").Value = 7
$ws.Range("C42").Value = 10
$ws.Range("C43").Value = 13
$ws.Range("C45").Value = 8
$ws.Range("C46").Value = 11
$ws.Range("C47").Value = 8
$ws.Range("C48").Value = 11
$ws.Range("C49").Value = 7
$ws.Range("C50").Value = 9
$ws.Range("C51").Value = 7
$ws.Range("C52").Value = 4
